# Apply updates to column F (dSF) on Sheet1 per commit: "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = 3
$ws.Range("F13").Value = 2
$ws.Range("F16").Value = 1
